$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D:K -> E:L), for the new FY2018 data.
$ws.Columns("D").EntireColumn.Insert()

# Copy the number formatting/styles from column E into the freshly inserted column D,
# one contiguous populated block at a time, so we match the look of the rest of the
# table (date format row 7, numeric format elsewhere) without materializing stray
# empty cells on rows that have no data (e.g. the blank spacer rows / section titles).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Populate the new column D with the latest (FY2018) figures ----
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 75700
$ws.Range("D9").Value = 64800
$ws.Range("D10").Value = 11000
$ws.Range("D12").Value = 10800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 96500
$ws.Range("D18").Value = -20700
$ws.Range("D20").Value = 1300
$ws.Range("D21").Value = -5000
$ws.Range("D22").Value = 300
$ws.Range("D23").Value = -19800
$ws.Range("D24").Value = -1300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -18500
$ws.Range("D27").Value = -18500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1300
$ws.Range("D33").Value = -19200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -19200

$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 11900
$ws.Range("D42").Value = 25500
$ws.Range("D43").Value = 18600
$ws.Range("D44").Value = 18800
$ws.Range("D45").Value = 1900
$ws.Range("D46").Value = 76700
$ws.Range("D47").Value = 4700
$ws.Range("D48").Value = 73200
$ws.Range("D49").Value = 12300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 32200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 199100
$ws.Range("D57").Value = 4100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 10600
$ws.Range("D60").Value = 14700
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 7800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 22500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 106000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 176600
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -19200
$ws.Range("D83").Value = 14400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -10400
$ws.Range("D91").Value = -8200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 7800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 100
$ws.Range("D101").Value = -600
$ws.Range("D102").Value = -3200

# Row 29 ("Effect Of Accounting Changes") historical years are marked as
# not-applicable rather than 0 in the refreshed data.
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
